$wb = $excel.ActiveWorkbook

# Sheet "Means" - Total Cancer Risk (row 9) and Total Respiratory (row 10)
$wsMeans = $wb.Worksheets.Item("Means")
$wsMeans.Range("B9").Value = 23
$wsMeans.Range("C9").Value = 18
$wsMeans.Range("B10").Value = 0.27
$wsMeans.Range("C10").Value = 0.2
$wsMeans.Range("G10").Value = 0.2

# Sheet "Standard Deviations" - Total Cancer Risk (row 9) and Total Respiratory (row 10)
$wsSD = $wb.Worksheets.Item("Standard Deviations")
$wsSD.Range("B9").Value = 7.2
$wsSD.Range("C9").Value = 4.2
$wsSD.Range("B10").Value = 0.094
$wsSD.Range("C10").Value = 0.024
$wsSD.Range("G10").Value = 0.000000000000000013
